$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.706082582473755
$ws.Range("B1").Value = 3.494296789169312
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 3.353888750076294
$ws.Range("E1").Value = 2.839561223983765
